# Final Universal code edit
# - Remove the first data row (old A2 = 15003), shifting all subsequent
#   facility counts up by one row (old A3..A53 -> new A2..A52).
# - Select the data range A1:A52.
# - Add a new, currently-empty cell at A55 with wrap-text formatting
#   applied (creates the second cellXfs entry used by that cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting row 2 shifts rows 3:53 up to 2:52, removing the 15003 entry
# and leaving the sheet with rows 1 (header) through 52 (value 36).
$ws.Rows(2).Delete()

# Restore/update the visible selection to the full (now-shifted) data range.
$ws.Range("A1:A52").Select()

# New blank cell two rows below the data, formatted with wrap text -
# this is what introduces the second <xf> in cellXfs.
$ws.Range("A55").WrapText = $true
